# Remove the "Solution Overview" slide (the one with the "✓ Semantic Search..."
# checkmark bullets), which sat at position 3 in the deck. Every slide that
# followed it shifts up by one position as a natural consequence of the
# deletion; no other slide content changes.
$p = $ppt.ActivePresentation
$p.Slides.Item(3).Delete()
